# Add Header row setting: insert a new blank row above row 1, pushing the
# existing header + data rows down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Insert() | Out-Null

# Move the selection cursor like the captured session (selection ends up on C8).
$ws.Range("C8").Select() | Out-Null
